# Model_stats.xlsx cleanup: reword the "Original model" description cells
# (row 5) to clarify the batch-input shape, and leave the selection on the
# cell that was last edited (B5), matching the author's final interactive
# state in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "1941 times there is a batch input, each batch has all 30,490 items included"
$ws.Range("C5").Value = "1941 batch inputs (1 for each day) á (30,490; 7; 20)?"

# Leave the active selection on B5 (bottom-right frozen pane), matching the
# saved workbook view in the commit.
$ws.Range("B5").Select()
